# Scoreboard.xlsx edit: remove the "Ledig plass" (placeholder) team entry
# from row 11 of the "ScoreM" list and the matching row 11 of the "Teams"
# lookup table, then leave the "Teams" sheet active/selected.

$wb = $excel.ActiveWorkbook

# --- 1) "ScoreM" sheet: delete the placeholder row (shifts rows 12-25 up) ---
$wsScoreM = $wb.Worksheets("ScoreM")
$wsScoreM.Rows(11).Delete()

# Leave the selection where the user ended up after the delete.
$wsScoreM.Range("C16").Select()

# --- 2) "Teams" sheet: the D:E lookup columns have the same placeholder   ---
#        row as ScoreM, but columns A:B hold an unrelated list that must   ---
#        stay put, so shift only D:E up manually instead of deleting the   ---
#        whole row.
$wsTeams = $wb.Worksheets("Teams")
for ($r = 11; $r -le 24; $r++) {
    $src = $r + 1
    $wsTeams.Range("D$r").Value = $wsTeams.Range("D$src").Value2
    $wsTeams.Range("E$r").Value = $wsTeams.Range("E$src").Value2
}
$wsTeams.Range("D25:E25").ClearContents()

# "Teams" ends up the active/selected sheet.
$wsTeams.Activate()
$wsTeams.Range("F9").Select()
